$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.979.48"
$ws.Range("E2").Value = "  -14.51%  "

# Row 3
$ws.Range("D3").Value = "2.316.02"
$ws.Range("E3").Value = "  -20.34%  "

# Row 4
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.23%  "

# Row 5
$ws.Range("D5").Value = "'434.67"
$ws.Range("E5").Value = "  -17.39%  "

# Row 6
$ws.Range("D6").Value = "'121.95"
$ws.Range("E6").Value = "  -15.46%  "

# Row 7
$ws.Range("D7").Value = "'0.994"
$ws.Range("E7").Value = "  -0.51%  "

# Row 8
$ws.Range("D8").Value = "'0.464"
$ws.Range("E8").Value = "  -15.09%  "

# Row 9
$ws.Range("D9").Value = "2.322.69"
$ws.Range("E9").Value = "  -20.33%  "

# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.0893"
$ws.Range("E10").Value = "  -16.78%  "

# Row 11
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "'5.09"
$ws.Range("E11").Value = "  -17.23%  "

# Row 12
$ws.Range("D12").Value = "'0.305"
$ws.Range("E12").Value = "  -14.91%  "

# Row 13
$ws.Range("D13").Value = "'0.120"
$ws.Range("E13").Value = "  -6.45%  "

# Row 14
$ws.Range("D14").Value = "2.668.77"
$ws.Range("E14").Value = "  -21.85%  "

# Row 15
$ws.Range("D15").Value = "51.762.05"
$ws.Range("E15").Value = "  -14.90%  "

# Row 16
$ws.Range("D16").Value = "'18.84"
$ws.Range("E16").Value = "  -16.45%  "

# Row 17
$ws.Range("D17").Value = "'0.0000119"
$ws.Range("E17").Value = "  -15.68%  "

# Row 18
$ws.Range("D18").Value = "2.293.89"
$ws.Range("E18").Value = "  -21.18%  "

# Row 19
$ws.Range("D19").Value = "'4.05"
$ws.Range("E19").Value = "  -17.45%  "

# Row 20
$ws.Range("D20").Value = "'296.37"
$ws.Range("E20").Value = "  -16.10%  "

# Row 21
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "'8.85"
$ws.Range("E21").Value = "  -23.35%  "

# Row 22
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.27%  "

# Row 23
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'5.65"
$ws.Range("E23").Value = "  -0.15%  "

# Row 24
$ws.Range("D24").Value = "'5.17"
$ws.Range("E24").Value = "  -20.55%  "

# Row 25
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'53.27"
$ws.Range("E25").Value = "  -17.84%  "

# Row 26
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'0.995"
$ws.Range("E26").Value = "  -0.21%  "

# Row 27
$ws.Range("D27").Value = "'0.371"
$ws.Range("E27").Value = "  -17.58%  "

# Row 28
$ws.Range("D28").Value = "2.329.14"
$ws.Range("E28").Value = "  -23.23%  "

# Row 29
$ws.Range("E29").Value = "  -25.12%  "

# Row 30
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.18%  "

# Row 31
$ws.Range("D31").Value = "'6.79"
$ws.Range("E31").Value = "  -12.90%  "

# Row 32
$ws.Range("D32").Value = "0.0₃0666"
$ws.Range("E32").Value = "  -23.05%  "

# Row 33
$ws.Range("D33").Value = "'143.86"
$ws.Range("E33").Value = "  -6.20%  "

# Row 34
$ws.Range("D34").Value = "'16.95"
$ws.Range("E34").Value = "  -13.74%  "

# Row 35
$ws.Range("D35").Value = "'1.32"
$ws.Range("E35").Value = "  -21.65%  "

# Row 36
$ws.Range("D36").Value = "'4.79"
$ws.Range("E36").Value = "  -14.23%  "

# Row 37
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'3.32"
$ws.Range("E37").Value = "  -24.51%  "

# Row 38
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.993"
$ws.Range("E38").Value = "  -17.21%  "

# Row 39
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'0.990"
$ws.Range("E39").Value = "  -0.67%  "

# Row 40
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "'0.762"
$ws.Range("E40").Value = "  -23.39%  "

# Row 41
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'31.92"
$ws.Range("E41").Value = "  -15.06%  "

# Row 42
$ws.Range("D42").Value = "'10.13"
$ws.Range("E42").Value = "  -2.02%  "

# Row 43
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'3.17"
$ws.Range("E43").Value = "  -14.22%  "

# Row 44
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.550"
$ws.Range("E44").Value = "  -15.74%  "

# Row 45
$ws.Range("D45").Value = "'0.0501"
$ws.Range("E45").Value = "  -14.08%  "

# Row 46
$ws.Range("D46").Value = "1.887.00"
$ws.Range("E46").Value = "  -17.63%  "

# Row 47
$ws.Range("D47").Value = "'1.15"
$ws.Range("E47").Value = "  -21.40%  "

# Row 48
$ws.Range("D48").Value = "'0.0203"
$ws.Range("E48").Value = "  -14.22%  "

# Row 49
$ws.Range("D49").Value = "'0.0809"
$ws.Range("E49").Value = "  -11.47%  "

# Row 50
$ws.Range("D50").Value = "'15.73"
$ws.Range("E50").Value = "  -22.76%  "

# Row 51
$ws.Range("D51").Value = "'3.95"
$ws.Range("E51").Value = "  -19.97%  "
